# Applies odds corrections to Sheet1 as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "Q2" = 1.7
    "R2" = 2.1
    "Q5" = 2.88
    "R5" = 1.4
    "BD11" = 151
    "N13" = 10
    "AA14" = 17
    "AO14" = 12
    "G14" = 2.2
    "I14" = 3.1
    "W14" = 8
    "BD15" = 151
    "Q15" = 1.9
    "R15" = 1.95
    "AB17" = 21
    "AC17" = 13
    "AE17" = 13
    "AG17" = 151
    "AH17" = 15
    "AM17" = 34
    "AN17" = 4
    "AP17" = 17
    "AR17" = 41
    "AS17" = 101
    "AU17" = 7.5
    "AX17" = 21
    "AY17" = 26
    "AZ17" = 67
    "BB17" = 151
    "H17" = 3.75
    "J17" = 2.3
    "L17" = 4.5
    "N17" = 13
    "O17" = 1.2
    "P17" = 4.33
    "Q17" = 1.7
    "R17" = 2.1
    "U17" = 1.67
    "V17" = 2.1
    "W17" = 8.5
    "X17" = 9.5
    "AP18" = 23
    "AT18" = 2.63
    "R18" = 1.75
    "S18" = 1.44
    "T18" = 2.63
    "AB20" = 29
    "AC20" = 12
    "AI20" = 11
    "AM20" = 23
    "AO20" = 17
    "AP20" = 23
    "AR20" = 67
    "G20" = 3.3
    "M20" = 1.05
    "N20" = 11
    "Q20" = 1.88
    "R20" = 1.98
    "Q21" = 1.67
    "R21" = 2.15
    "N22" = 8
    "Q22" = 2.25
    "R22" = 1.62
    "S22" = 1.5
    "S23" = 1.25
    "AC24" = 9
    "AD24" = 7.5
    "AH24" = 11
    "AN24" = 3.5
    "AO24" = 9
    "AQ24" = 29
    "G24" = 1.7
    "H24" = 3.9
    "J24" = 2.38
    "M24" = 1.07
    "N24" = 9
    "S24" = 1.44
    "T24" = 2.63
    "Y24" = 8.5
    "S25" = 1.36
    "M26" = 1.07
    "N26" = 9
    "O26" = 1.36
    "P26" = 3
    "Q26" = 2.1
    "R26" = 1.7
    "O27" = 1.3
    "P27" = 3.4
    "Q27" = 2.03
    "R27" = 1.83
    "Q56" = 1.8
    "R56" = 2
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value2 = $updates[$cellRef]
}
